# Apply the marks-sheet corrections described in the commit:
#   "Some further error corrections when working with the marks sheet"
#
# 1. Cover page: TERM I -> TERM II
# 2. Marks table: GEOGRAPHY and AGRICULTURE rows swap places (including
#    a couple of value corrections), ENTREPRENEURSHIP's marks shift from
#    its Paper 1 row down to its Paper 2 row (with corrections), SUBMATH's
#    row is blanked out, and GENERAL PAPER / TOTAL POINTS get corrected
#    values.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Cover page title: "TERM I" -> "TERM II"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("TERM I", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TERM II", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Marks table edits (Tables(1))
# ---------------------------------------------------------------------
$t = $d.Tables(1)

# --- Subject block 1 (rows 4-6): was GEOGRAPHY, becomes AGRICULTURE,
#     keeping only a (corrected) "Paper 2" entry.
$t.Cell(4, 1).Range.Text = "AGRICULTURE"   # Subject name
$t.Cell(4, 2).Range.Text = ""              # Paper 1 -> blank
$t.Cell(4, 4).Range.Text = ""              # Total marks -> blank
$t.Cell(4, 6).Range.Text = ""              # Paper grade -> blank
# row4 col8 (Subject grade, 'O') stays unchanged

$t.Cell(5, 4).Range.Text = "65.0"          # Paper 2 total marks 61 -> 65.0
$t.Cell(5, 6).Range.Text = "C5"            # Paper 2 grade C6 -> C5
# row5 col2 "Paper 2" label stays unchanged

$t.Cell(6, 2).Range.Text = ""              # Paper 3 -> blank
$t.Cell(6, 4).Range.Text = ""              # Total marks -> blank
$t.Cell(6, 6).Range.Text = ""              # Paper grade -> blank

# --- Subject block 2 (rows 7-9): ENTREPRENEURSHIP -> ENTERPRENEURSHIP
#     (spelling fix); its marks move from the Paper 1 row to the Paper 2 row.
$t.Cell(7, 1).Range.Text = "ENTERPRENEURSHIP"
$t.Cell(7, 2).Range.Text = ""              # Paper 1 -> blank
$t.Cell(7, 4).Range.Text = ""              # Total marks -> blank
$t.Cell(7, 6).Range.Text = ""              # Paper grade -> blank
$t.Cell(7, 8).Range.Text = "O"             # Subject grade F9 -> O

$t.Cell(8, 2).Range.Text = "Paper 2"       # Paper label blank -> Paper 2
$t.Cell(8, 4).Range.Text = "45.0"          # Total marks blank -> 45.0
$t.Cell(8, 6).Range.Text = "P8"            # Paper grade blank -> P8

# row9 (Paper 3 line) stays fully blank

# --- Subject block 3 (rows 10-12): was AGRICULTURE, becomes GEOGRAPHY,
#     all paper rows blanked except the subject grade correction.
$t.Cell(10, 1).Range.Text = "GEOGRAPHY"
$t.Cell(10, 2).Range.Text = ""             # Paper 1 -> blank
$t.Cell(10, 4).Range.Text = ""             # Total marks -> blank
$t.Cell(10, 6).Range.Text = ""             # Paper grade -> blank
$t.Cell(10, 8).Range.Text = "X"            # Subject grade C5 -> X

$t.Cell(11, 2).Range.Text = "Paper 2"      # Paper label blank -> Paper 2
$t.Cell(11, 4).Range.Text = ""             # Total marks 67 -> blank
$t.Cell(11, 6).Range.Text = ""             # Paper grade 67 -> blank

$t.Cell(12, 2).Range.Text = "Paper 3"      # Paper label blank -> Paper 3
$t.Cell(12, 4).Range.Text = ""             # Total marks 67 -> blank
$t.Cell(12, 6).Range.Text = ""             # Paper grade 67 -> blank

# --- Subject block 4 (row 13): SUBMATH row blanked out entirely.
$t.Cell(13, 1).Range.Text = ""             # Subject name -> blank
$t.Cell(13, 4).Range.Text = ""             # Total marks -> blank
$t.Cell(13, 6).Range.Text = ""             # Paper grade -> blank
$t.Cell(13, 8).Range.Text = ""             # Subject grade -> blank
# row13 col2 "Paper 1" label stays unchanged

# --- Subject block 5 (row 14): GENERAL PAPER values corrected.
$t.Cell(14, 4).Range.Text = "30.0"         # Total marks 42 -> 30.0
$t.Cell(14, 6).Range.Text = "F9"           # Paper grade P8 -> F9
$t.Cell(14, 8).Range.Text = "F9"           # Subject grade P8 -> F9

# --- TOTAL POINTS row (row 15): value corrected.
$t.Cell(15, 4).Range.Text = "X"            # Total points 2 -> X
